$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = "Fallo"
$ws.Range("H9").Value = -1

$ws.Range("G13").Value = "Fallo"
$ws.Range("H13").Value = -1

$ws.Range("G25").Value = "Fallo"
$ws.Range("H25").Value = -1

$ws.Range("A28").Value = 14631199
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "2025-09-07"
$ws.Range("C28").Value = "Hanna Chang"
$ws.Range("D28").Value = "Meiqi Guo"
$ws.Range("E28").Value = "Gana Hanna Chang"
$ws.Range("F28").Value = 1.83
